$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.359.78'
$ws.Range("E2").Value = '  -3.00%  '
$ws.Range("D3").Value = '2.286.32'
$ws.Range("E3").Value = '  -2.97%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '494.46'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.32'
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("D9").Value = '2.287.17'
$ws.Range("E9").Value = '  -3.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0942'
$ws.Range("E10").Value = '  -3.42%  '
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("E13").Value = '  -3.75%  '
$ws.Range("D14").Value = '2.692.82'
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.53'
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = '54.198.72'
$ws.Range("E16").Value = '  -3.19%  '
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("D18").Value = '2.283.55'
$ws.Range("E18").Value = '  -3.29%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.05'
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '302.60'
$ws.Range("E21").Value = '  -2.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.41'
$ws.Range("E22").Value = '  +3.21%  '
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.39'
$ws.Range("E24").Value = '  -2.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.68'
$ws.Range("E25").Value = '  -2.63%  '
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("D28").Value = '2.394.14'
$ws.Range("E28").Value = '  -2.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.148'
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.09'
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.84'
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("D33").Value = '0.0₃0683'
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.87'
$ws.Range("E34").Value = '  +1.78%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.08'
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("E39").Value = '  +0.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.870'
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.62'
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '126.56'
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.81'
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0889'
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '238.61'
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("E51").Value = '  +0.05%  '
